# Adds three new "Run 1/Run 2/Run 3" queue-size-10 sample columns to both
# summary tables on Sheet1, matching commit "Added more tests for varying queue size".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: the old "Keyed" sub-label in C2 is removed (blank cell, style kept) ---
$ws.Range("C2").Value = ""

# --- Row 4: new Run 1/Run 2/Run 3 sub-headers (D4:F4 already carry the right style) ---
$ws.Range("D4").Value = "Run 1"
$ws.Range("E4").Value = "Run 2"
$ws.Range("F4").Value = "Run 3"

# --- Shift columns E:J two places right to G:L for the first summary table (rows 5-9) ---
# (right-to-left per column so we never clobber a not-yet-moved source cell)
$ws.Range("J5").Copy($ws.Range("L5"))
$ws.Range("I5").Copy($ws.Range("K5"))
$ws.Range("H5").Copy($ws.Range("J5"))
$ws.Range("G5").Copy($ws.Range("I5"))
$ws.Range("F5").Copy($ws.Range("H5"))
$ws.Range("E5").Copy($ws.Range("G5"))

$ws.Range("J6").Copy($ws.Range("L6"))
$ws.Range("I6").Copy($ws.Range("K6"))
$ws.Range("H6").Copy($ws.Range("J6"))
$ws.Range("G6").Copy($ws.Range("I6"))
$ws.Range("F6").Copy($ws.Range("H6"))
$ws.Range("E6").Copy($ws.Range("G6"))

$ws.Range("J7").Copy($ws.Range("L7"))
$ws.Range("I7").Copy($ws.Range("K7"))
$ws.Range("H7").Copy($ws.Range("J7"))
$ws.Range("G7").Copy($ws.Range("I7"))
$ws.Range("F7").Copy($ws.Range("H7"))
$ws.Range("E7").Copy($ws.Range("G7"))

$ws.Range("J8").Copy($ws.Range("L8"))
$ws.Range("I8").Copy($ws.Range("K8"))
$ws.Range("H8").Copy($ws.Range("J8"))
$ws.Range("G8").Copy($ws.Range("I8"))
$ws.Range("F8").Copy($ws.Range("H8"))
$ws.Range("E8").Copy($ws.Range("G8"))

$ws.Range("J9").Copy($ws.Range("L9"))
$ws.Range("I9").Copy($ws.Range("K9"))
$ws.Range("H9").Copy($ws.Range("J9"))
$ws.Range("G9").Copy($ws.Range("I9"))
$ws.Range("F9").Copy($ws.Range("H9"))
$ws.Range("E9").Copy($ws.Range("G9"))

# --- Shift columns E:J two places right to G:L for the second summary table (rows 16-19) ---
# (right-to-left per column so we never clobber a not-yet-moved source cell)
$ws.Range("J16").Copy($ws.Range("L16"))
$ws.Range("I16").Copy($ws.Range("K16"))
$ws.Range("H16").Copy($ws.Range("J16"))
$ws.Range("G16").Copy($ws.Range("I16"))
$ws.Range("F16").Copy($ws.Range("H16"))
$ws.Range("E16").Copy($ws.Range("G16"))

$ws.Range("J17").Copy($ws.Range("L17"))
$ws.Range("I17").Copy($ws.Range("K17"))
$ws.Range("H17").Copy($ws.Range("J17"))
$ws.Range("G17").Copy($ws.Range("I17"))
$ws.Range("F17").Copy($ws.Range("H17"))
$ws.Range("E17").Copy($ws.Range("G17"))

$ws.Range("J18").Copy($ws.Range("L18"))
$ws.Range("I18").Copy($ws.Range("K18"))
$ws.Range("H18").Copy($ws.Range("J18"))
$ws.Range("G18").Copy($ws.Range("I18"))
$ws.Range("F18").Copy($ws.Range("H18"))
$ws.Range("E18").Copy($ws.Range("G18"))

$ws.Range("J19").Copy($ws.Range("L19"))
$ws.Range("I19").Copy($ws.Range("K19"))
$ws.Range("H19").Copy($ws.Range("J19"))
$ws.Range("G19").Copy($ws.Range("I19"))
$ws.Range("F19").Copy($ws.Range("H19"))
$ws.Range("E19").Copy($ws.Range("G19"))

# --- New E/F columns take on column C/D's formatting, then get their own values ---
$ws.Range("C5").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 10

$ws.Range("C6").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("E6").Value = 1000
$ws.Range("F6").Value = 1000

$ws.Range("C7").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("E7").Value = 415
$ws.Range("F7").Value = 387

$ws.Range("C8").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 11

$ws.Range("C9").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("E9").Value = 60
$ws.Range("F9").Value = 64

$ws.Range("C16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 10

$ws.Range("C17").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("D17").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 5

$ws.Range("C18").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("D18").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("E18").Value = 2425571
$ws.Range("F18").Value = 3040513

$ws.Range("C19").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("D19").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$ws.Range("E19").Value = 74
$ws.Range("F19").Value = 75

# --- Row 10: footnote row. "10 hours 24 mins" note moves from E10 to G10; ---
# --- E10/F10 are vacated entirely, and K10/L10 gain blank styled cells to match ---
# --- the table's new width. ---
$ws.Range("J10").Copy($ws.Range("K10"))
$ws.Range("J10").Copy($ws.Range("L10"))
$ws.Range("E10").Copy($ws.Range("G10"))
$ws.Range("E10").Clear()
$ws.Range("F10").Clear()

# --- New row 15: Run 1/Run 2/Run 3 labels above the second table's data rows ---
$ws.Range("D15").Value = "Run 1"
$ws.Range("E15").Value = "Run 2"
$ws.Range("F15").Value = "Run 3"

# --- Restore the active selection Excel leaves behind after this edit ---
[void]$ws.Range("F3").Select()

